$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.156.70"
$ws.Range("E2").Value = "  -3.21%  "

$ws.Range("D3").Value = "1.608.82"
$ws.Range("E3").Value = "  -2.54%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9997"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "301.56"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.39%  "

$ws.Range("E7").Value = "  -2.83%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3655"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "49.16"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -5.29%  "

$ws.Range("E10").Value = "  +0.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.270"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -6.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08083"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.09%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.07"
$ws.Range("D13").ClearFormats()

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.611"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -6.68%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.453"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -6.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001254"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.82%  "

$ws.Range("D17").Value = "1.613.77"
$ws.Range("E17").Value = "  -2.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.52"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.40%  "

$ws.Range("E19").Value = "  -2.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.40"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -6.72%  "

$ws.Range("E21").Value = "  -5.23%  "

$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.10"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.61%  "

$ws.Range("D24").Value = "23.184.61"
$ws.Range("E24").Value = "  -3.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.352"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.869"
$ws.Range("D26").ClearFormats()

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.11"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.59"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.274"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.11"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.415"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.874"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -12.76%  "

$ws.Range("D33").Value = "1.790.82"
$ws.Range("E33").Value = "  -1.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9717"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -7.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07710"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02771"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -6.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2565"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.36%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.252"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -6.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.15"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -7.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08904"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.16%  "

$ws.Range("E41").Value = "  -2.42%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7209"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.84"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.76"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6685"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.314"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -6.25%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9994"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.25%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.979"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.41%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08016"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "130.97"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.67%  "

$ws.Range("E51").Value = "  -2.95%  "
